# Applies the commit "ConsoleIO now has something to check for commas":
#   1. Inserts a new paragraph "Did something with ConsoleIO to correct
#      comma input" (with "Did something with ConsoleIO" highlighted
#      cyan) right before the "Create models for both Reports..."
#      paragraph. The page-break marker that used to render at the start
#      of "Create models..." now renders at the start of this new
#      paragraph instead.
#   2. Because the document got one paragraph longer, the page-break
#      marker that used to fall on "All financial math must use decimal."
#      now falls one paragraph earlier, on "Strive to generate reports
#      with LINQ...".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: new "Did something with ConsoleIO..." paragraph
# ---------------------------------------------------------------------
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Create models for both Reports*") {
        $targetIdx = $i
        break
    }
}

$createModelsPara = $d.Paragraphs.Item($targetIdx)
$createModelsRange = $createModelsPara.Range

$newParaXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="464A1D98" w14:textId="2AADD5EE" w:rsidR="00F01503" w:rsidRDefault="00F01503" w:rsidP="00AF530D"><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:lastRenderedPageBreak/><w:t>D</w:t></w:r><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>id something with ConsoleIO</w:t></w:r><w:r><w:t xml:space="preserve"> to correct comma input</w:t></w:r></w:p>
<w:p w14:paraId="464A1D99" w14:textId="2AADD54E" w:rsidR="00F01503" w:rsidRDefault="00F01503" w:rsidP="00AF530D"><w:r><w:t>Create models for both Reports (CategoryValue &amp; )</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$createModelsRange.InsertXML($newParaXml) | Out-Null

# ---------------------------------------------------------------------
# Change 2: move the lastRenderedPageBreak marker up one paragraph,
# from "All financial math..." to "Strive to generate reports..."
# ---------------------------------------------------------------------
$striveIdx = -1
$financialIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*Strive to generate reports with LINQ*") {
        $striveIdx = $i
    }
    if ($text -like "*All financial math*") {
        $financialIdx = $i
    }
}

# Remove the marker from "All financial math must use decimal."
$financialPara = $d.Paragraphs.Item($financialIdx)
$financialRange = $financialPara.Range
$financialXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="389A0B95" w14:textId="77777777" w:rsidR="00C16017" w:rsidRPr="00C16017" w:rsidRDefault="00C16017" w:rsidP="00C16017"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r w:rsidRPr="00C16017"><w:t>All financial math must use&#160;decimal.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$financialRange.InsertXML($financialXml) | Out-Null

# Add the marker to "Strive to generate reports with LINQ..."
$strivePara = $d.Paragraphs.Item($striveIdx)
$striveRange = $strivePara.Range
$striveXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="2AFCA030" w14:textId="77777777" w:rsidR="00C16017" w:rsidRPr="00C16017" w:rsidRDefault="00C16017" w:rsidP="00C16017"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r w:rsidRPr="00C16017"><w:lastRenderedPageBreak/><w:t>Strive to generate reports with LINQ. If you run into too much friction, solve the problem with loops and intermediate collections.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$striveRange.InsertXML($striveXml) | Out-Null
